$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item 31 (row 36): "cable lug 0,75mm^2" -> now covers two small + four medium
# lugs, quantity bumped from 2 to 6.
$ws.Range("C36").Value = 6
$ws.Range("D36").Value = "cable lug 0,75mm^2, two small, four medium"

# Bottom total-row label: "SUMME" -> "total"
$ws.Range("D39").Value = "total"

# Cosmetic view refresh: re-centre/zoom the sheet and move the selection.
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 9
$win.ScrollColumn = 1
[void]$ws.Range("D40").Select()
